# Insert a new data row at row 91 (pushing existing rows 91:177 down to 92:178)
# and populate it with a new Kiwi price record, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 91; this shifts rows 91-177 down to 92-178
# and copies formatting from the row above (keeps the date-format style on column D).
$ws.Rows("91:91").Insert()

# Populate the newly inserted row 91 with the new record's data.
$ws.Range("A91").Value = 5
$ws.Range("B91").Value = "Macroferia Regional de Talca"
$ws.Range("C91").Value = "Maule"
$ws.Range("D91").Value = 44452
$ws.Range("E91").Value = 7
$ws.Range("F91").Value = "Fruta"
$ws.Range("G91").Value = 100101
$ws.Range("H91").Value = "Berries"
$ws.Range("I91").Value = 100101007
$ws.Range("J91").Value = "Kiwi"
$ws.Range("K91").Value = "Hayward"
$ws.Range("L91").Value = "Primera"
$ws.Range("M91").Value = 300
$ws.Range("N91").Value = 12000
$ws.Range("O91").Value = 12000
$ws.Range("P91").Value = 12000
$ws.Range("Q91").Value = "$/bandeja 18 kilos"
$ws.Range("R91").Value = "Provincia de Curicó"
$ws.Range("S91").Value = 667
$ws.Range("T91").Value = 18
